$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values as per diff: row 10 col A, row 12 col A, row 13 col D, row 18 col A
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("A18").Value = -21.694
